# This script reproduces the "Luan them phan doi cua cung vo chinh dieu"
# (Add opposite-palace interpretation for palaces with no main star) update:
# it appends 105 new "tai cung doi Tat Ach" rows/strings to the Sheet2 table,
# and updates the sheet selection to the end of the newly written range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Tử Vi tọa thủ tại cung đối Tật Ách',
    'Thiên Cơ tọa thủ tại cung đối Tật Ách',
    'Thái Dương tọa thủ tại cung đối Tật Ách',
    'Vũ Khúc tọa thủ tại cung đối Tật Ách',
    'Thiên Đồng tọa thủ tại cung đối Tật Ách',
    'Liêm Trinh tọa thủ tại cung đối Tật Ách',
    'Thiên Phủ tọa thủ tại cung đối Tật Ách',
    'Thái Âm tọa thủ tại cung đối Tật Ách',
    'Tham Lang tọa thủ tại cung đối Tật Ách',
    'Cự Môn tọa thủ tại cung đối Tật Ách',
    'Thiên Tướng tọa thủ tại cung đối Tật Ách',
    'Thiên Lương tọa thủ tại cung đối Tật Ách',
    'Thất Sát tọa thủ tại cung đối Tật Ách',
    'Phá Quân tọa thủ tại cung đối Tật Ách',
    'Tử Vi đồng cung Thiên Cơ tại cung đối Tật Ách',
    'Tử Vi đồng cung Thái Dương tại cung đối Tật Ách',
    'Tử Vi đồng cung Vũ Khúc tại cung đối Tật Ách',
    'Tử Vi đồng cung Thiên Đồng tại cung đối Tật Ách',
    'Tử Vi đồng cung Liêm Trinh tại cung đối Tật Ách',
    'Tử Vi đồng cung Thiên Phủ tại cung đối Tật Ách',
    'Tử Vi đồng cung Thái Âm tại cung đối Tật Ách',
    'Tử Vi đồng cung Tham Lang tại cung đối Tật Ách',
    'Tử Vi đồng cung Cự Môn tại cung đối Tật Ách',
    'Tử Vi đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Tử Vi đồng cung Thiên Lương tại cung đối Tật Ách',
    'Tử Vi đồng cung Thất Sát tại cung đối Tật Ách',
    'Tử Vi đồng cung Phá Quân tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thái Dương tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Vũ Khúc tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thiên Đồng tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Liêm Trinh tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thiên Phủ tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thái Âm tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Tham Lang tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Cự Môn tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thiên Lương tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Thất Sát tại cung đối Tật Ách',
    'Thiên Cơ đồng cung Phá Quân tại cung đối Tật Ách',
    'Thái Dương đồng cung Vũ Khúc tại cung đối Tật Ách',
    'Thái Dương đồng cung Thiên Đồng tại cung đối Tật Ách',
    'Thái Dương đồng cung Liêm Trinh tại cung đối Tật Ách',
    'Thái Dương đồng cung Thiên Phủ tại cung đối Tật Ách',
    'Thái Dương đồng cung Thái Âm tại cung đối Tật Ách',
    'Thái Dương đồng cung Tham Lang tại cung đối Tật Ách',
    'Thái Dương đồng cung Cự Môn tại cung đối Tật Ách',
    'Thái Dương đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Thái Dương đồng cung Thiên Lương tại cung đối Tật Ách',
    'Thái Dương đồng cung Thất Sát tại cung đối Tật Ách',
    'Thái Dương đồng cung Phá Quân tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Thiên Đồng tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Liêm Trinh tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Thiên Phủ tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Thái Âm tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Tham Lang tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Cự Môn tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Thiên Lương tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Thất Sát tại cung đối Tật Ách',
    'Vũ Khúc đồng cung Phá Quân tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Liêm Trinh tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Thiên Phủ tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Thái Âm tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Tham Lang tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Cự Môn tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Thiên Lương tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Thất Sát tại cung đối Tật Ách',
    'Thiên Đồng đồng cung Phá Quân tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Thiên Phủ tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Thái Âm tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Tham Lang tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Cự Môn tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Thiên Lương tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Thất Sát tại cung đối Tật Ách',
    'Liêm Trinh đồng cung Phá Quân tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Thái Âm tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Tham Lang tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Cự Môn tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Thiên Lương tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Thất Sát tại cung đối Tật Ách',
    'Thiên Phủ đồng cung Phá Quân tại cung đối Tật Ách',
    'Thái Âm đồng cung Tham Lang tại cung đối Tật Ách',
    'Thái Âm đồng cung Cự Môn tại cung đối Tật Ách',
    'Thái Âm đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Thái Âm đồng cung Thiên Lương tại cung đối Tật Ách',
    'Thái Âm đồng cung Thất Sát tại cung đối Tật Ách',
    'Thái Âm đồng cung Phá Quân tại cung đối Tật Ách',
    'Tham Lang đồng cung Cự Môn tại cung đối Tật Ách',
    'Tham Lang đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Tham Lang đồng cung Thiên Lương tại cung đối Tật Ách',
    'Tham Lang đồng cung Thất Sát tại cung đối Tật Ách',
    'Tham Lang đồng cung Phá Quân tại cung đối Tật Ách',
    'Cự Môn đồng cung Thiên Tướng tại cung đối Tật Ách',
    'Cự Môn đồng cung Thiên Lương tại cung đối Tật Ách',
    'Cự Môn đồng cung Thất Sát tại cung đối Tật Ách',
    'Cự Môn đồng cung Phá Quân tại cung đối Tật Ách',
    'Thiên Tướng đồng cung Thiên Lương tại cung đối Tật Ách',
    'Thiên Tướng đồng cung Thất Sát tại cung đối Tật Ách',
    'Thiên Tướng đồng cung Phá Quân tại cung đối Tật Ách',
    'Thiên Lương đồng cung Thất Sát tại cung đối Tật Ách',
    'Thiên Lương đồng cung Phá Quân tại cung đối Tật Ách',
    'Thất Sát đồng cung Phá Quân tại cung đối Tật Ách'
)

$startRow = 4331
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $text = $values[$i]
    $ws.Cells.Item($r, 1).Value = $text
    $ws.Cells.Item($r, 2).Value = $text
}

$lastRow = $startRow + $values.Count - 1

$ws.Activate()
$ws.Range("B4345:B" + $lastRow).Select()
